$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12
$ws.Range('A12').Value = 'ddl_011'
$ws.Range('B12').Value = 'y'
$ws.Range('C12').Value = '清空带有标量和向量索引的单分区空表'
$ws.Range('D12').Value = 'DDL'
$ws.Range('E12').Value = 'truncate'
$ws.Range('F12').Value = 'mixindex015'
$ws.Range('G12').NumberFormat = '@'
$ws.Range('H12').Value = 'truncate table $mixindex015'
$ws.Range('I12').Value = 'select /*+ vector_pre */ id,feature_id1,feature1_index$distance from vector($mixindex015, feature1, array[1.2949047088623047, 0.8476371169090271, 0.41010794043540955, 0.9785786271095276, 0.6064758896827698, 0.5890575647354126, 0.4460025131702423, 0.43684643507003784], 10) where name='''' order by feature1_index$distance limit 10'
$ws.Range('J12').Value = 'src/test/resources/io.dingodb.test/testdata/mysqlcases/ddl/expectedresult/ddl_011.csv'
$ws.Range('K12').Value = 'csv_equals'

# Row 13
$ws.Range('A13').Value = 'ddl_012'
$ws.Range('B13').Value = 'y'
$ws.Range('C13').Value = '清空带有标量和向量索引的单分区非空表'
$ws.Range('D13').Value = 'DDL'
$ws.Range('E13').Value = 'truncate'
$ws.Range('F13').Value = 'mixindex015'
$ws.Range('G13').Value = 'mix015_value1'
$ws.Range('H13').Value = 'truncate table $mixindex015'
$ws.Range('I13').Value = 'select id,name,feature_id2,feature2_index$distance from vector($mixindex015, feature2, array[1.3893274068832397, 0.4223838150501251, 0.7195155024528503, 0.28465819358825684, 0.9267012476921082, 0.4375186264514923, 0.7362583875656128, 0.41960853338241577, 0.7168405055999756, 0.5768887996673584, 0.4182721674442291, 0.129996195435524, 0.21704305708408356, 0.35628095269203186, 0.5203919410705566, 0.4046420454978943, 0.1657610833644867, 0.7787348031997681, 0.16017264127731323, 0.7188393473625183, 0.916609525680542, 0.6900423765182495, 0.6827380657196045, 0.49170464277267456, 0.5505375266075134, 0.33907604217529297, 0.09597073495388031, 0.008679530583322048, 0.9524646997451782, 0.2128734141588211, 0.6533687710762024, 0.6094813346862793, 0.18400055170059204, 0.9364618062973022, 0.7497748732566833, 0.05966084823012352, 0.40710607171058655, 0.6920192241668701, 0.649844765663147, 0.705480694770813, 0.11741353571414948, 0.2164693921804428, 0.48223137855529785, 0.6425648331642151, 0.5236963033676147, 0.6490180492401123, 0.5132198333740234, 0.5319958925247192, 0.5141375660896301, 0.9233165979385376, 0.5702359080314636, 0.46916520595550537, 0.7964460849761963, 0.17150050401687622, 0.4680892825126648, 0.6861740946769714, 0.4802965819835663, 0.9073042273521423, 0.17548426985740662, 0.9943628311157227, 0.9149019718170166, 0.7661579847335815, 0.13299474120140076, 0.7843778133392334], 10, map[efSearch, 40])'
$ws.Range('J13').Value = 'src/test/resources/io.dingodb.test/testdata/mysqlcases/ddl/expectedresult/ddl_012.csv'
$ws.Range('K13').Value = 'csv_equals'

# Row 14
$ws.Range('A14').Value = 'ddl_013'
$ws.Range('B14').Value = 'y'
$ws.Range('C14').Value = '清空带有标量和向量索引的多分区空表'
$ws.Range('D14').Value = 'DDL'
$ws.Range('E14').Value = 'truncate'
$ws.Range('F14').Value = 'mixindex016'
$ws.Range('G14').NumberFormat = '@'
$ws.Range('H14').Value = 'truncate table $mixindex016'
$ws.Range('I14').Value = 'select id,feature_id,feature_index$distance from vector($mixindex016, feature, array[1.3893274068832397, 0.4223838150501251, 0.7195155024528503, 0.28465819358825684, 0.9267012476921082, 0.4375186264514923, 0.7362583875656128, 0.41960853338241577, 0.7168405055999756, 0.5768887996673584, 0.4182721674442291, 0.129996195435524, 0.21704305708408356, 0.35628095269203186, 0.5203919410705566, 0.4046420454978943, 0.1657610833644867, 0.7787348031997681, 0.16017264127731323, 0.7188393473625183, 0.916609525680542, 0.6900423765182495, 0.6827380657196045, 0.49170464277267456, 0.5505375266075134, 0.33907604217529297, 0.09597073495388031, 0.008679530583322048, 0.9524646997451782, 0.2128734141588211, 0.6533687710762024, 0.6094813346862793, 0.18400055170059204, 0.9364618062973022, 0.7497748732566833, 0.05966084823012352, 0.40710607171058655, 0.6920192241668701, 0.649844765663147, 0.705480694770813, 0.11741353571414948, 0.2164693921804428, 0.48223137855529785, 0.6425648331642151, 0.5236963033676147, 0.6490180492401123, 0.5132198333740234, 0.5319958925247192, 0.5141375660896301, 0.9233165979385376, 0.5702359080314636, 0.46916520595550537, 0.7964460849761963, 0.17150050401687622, 0.4680892825126648, 0.6861740946769714, 0.4802965819835663, 0.9073042273521423, 0.17548426985740662, 0.9943628311157227, 0.9149019718170166, 0.7661579847335815, 0.13299474120140076, 0.7843778133392334], 10, map[efSearch, 40]) order by feature_index$distance limit 10'
$ws.Range('J14').Value = 'src/test/resources/io.dingodb.test/testdata/mysqlcases/ddl/expectedresult/ddl_013.csv'
$ws.Range('K14').Value = 'csv_equals'

# Row 15
$ws.Range('A15').Value = 'ddl_014'
$ws.Range('B15').Value = 'y'
$ws.Range('C15').Value = '清空带有标量和向量索引的多分区非空表'
$ws.Range('D15').Value = 'DDL'
$ws.Range('E15').Value = 'truncate'
$ws.Range('F15').Value = 'mixindex017'
$ws.Range('G15').Value = 'mix_common_value1'
$ws.Range('H15').Value = 'truncate table $mixindex017'
$ws.Range('I15').Value = 'select /*+ vector_pre */ id,birthday,feature_id,feature_index$distance from vector($mixindex017, feature, array[1.3893274068832397, 0.4223838150501251, 0.7195155024528503, 0.28465819358825684, 0.9267012476921082, 0.4375186264514923, 0.7362583875656128, 0.41960853338241577, 0.7168405055999756, 0.5768887996673584, 0.4182721674442291, 0.129996195435524, 0.21704305708408356, 0.35628095269203186, 0.5203919410705566, 0.4046420454978943, 0.1657610833644867, 0.7787348031997681, 0.16017264127731323, 0.7188393473625183, 0.916609525680542, 0.6900423765182495, 0.6827380657196045, 0.49170464277267456, 0.5505375266075134, 0.33907604217529297, 0.09597073495388031, 0.008679530583322048, 0.9524646997451782, 0.2128734141588211, 0.6533687710762024, 0.6094813346862793, 0.18400055170059204, 0.9364618062973022, 0.7497748732566833, 0.05966084823012352, 0.40710607171058655, 0.6920192241668701, 0.649844765663147, 0.705480694770813, 0.11741353571414948, 0.2164693921804428, 0.48223137855529785, 0.6425648331642151, 0.5236963033676147, 0.6490180492401123, 0.5132198333740234, 0.5319958925247192, 0.5141375660896301, 0.9233165979385376, 0.5702359080314636, 0.46916520595550537, 0.7964460849761963, 0.17150050401687622, 0.4680892825126648, 0.6861740946769714, 0.4802965819835663, 0.9073042273521423, 0.17548426985740662, 0.9943628311157227, 0.9149019718170166, 0.7661579847335815, 0.13299474120140076, 0.7843778133392334], 10, map[efSearch, 40]) where birthday=''2004-12-23'' order by feature_index$distance limit 10'
$ws.Range('J15').Value = 'src/test/resources/io.dingodb.test/testdata/mysqlcases/ddl/expectedresult/ddl_014.csv'
$ws.Range('K15').Value = 'csv_equals'

# Row 16
$ws.Range('A16').Value = 'ddl_015'
$ws.Range('B16').Value = 'y'
$ws.Range('C16').Value = '清空带有标量和向量索引的单分区非空表-标量查询'
$ws.Range('D16').Value = 'DDL'
$ws.Range('E16').Value = 'truncate'
$ws.Range('F16').Value = 'mixindex015'
$ws.Range('G16').Value = 'mix015_value1'
$ws.Range('H16').Value = 'truncate table $mixindex015'
$ws.Range('I16').Value = 'select id,name from $mixindex015 where name='''''
$ws.Range('J16').Value = 'src/test/resources/io.dingodb.test/testdata/mysqlcases/ddl/expectedresult/ddl_015.csv'
$ws.Range('K16').Value = 'csv_equals'

$ws.Range('I28').Select()
